# Weekly update: two new price records for the current week are inserted
# at the top of the historical table (rows 507-508), pushing every
# existing record down by two rows (507->509 ... 530->532).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 507:508 (existing rows 507-530 shift to 509-532).
$ws.Rows("507:508").Insert()

# New row 507
$ws.Range("A507").Value = 10
$ws.Range("B507").Value = "Vega Modelo de Temuco"
$ws.Range("C507").Value = "La Araucanía"
$ws.Range("D507").Value = "2022-07-11"
$ws.Range("E507").Value = 9
$ws.Range("F507").Value = 100112043
$ws.Range("G507").Value = "Pepino ensalada"
$ws.Range("H507").Value = "Sin especificar"
$ws.Range("I507").Value = "Primera"
$ws.Range("J507").Value = 270
$ws.Range("K507").Value = 20000
$ws.Range("L507").Value = 21000
$ws.Range("M507").Value = 20556
$ws.Range("N507").Value = "$/caja 50 unidades"
$ws.Range("O507").Value = "Región de Arica y Parinacota"
$ws.Range("P507").Value = 411
$ws.Range("Q507").Value = 50
$ws.Range("R507").Value = "Hortaliza"

# New row 508
$ws.Range("A508").Value = 10
$ws.Range("B508").Value = "Vega Modelo de Temuco"
$ws.Range("C508").Value = "La Araucanía"
$ws.Range("D508").Value = "2022-07-11"
$ws.Range("E508").Value = 9
$ws.Range("F508").Value = 100112043
$ws.Range("G508").Value = "Pepino ensalada"
$ws.Range("H508").Value = "Sin especificar"
$ws.Range("I508").Value = "Segunda"
$ws.Range("J508").Value = 60
$ws.Range("K508").Value = 18000
$ws.Range("L508").Value = 18000
$ws.Range("M508").Value = 18000
$ws.Range("N508").Value = "$/caja 80 unidades"
$ws.Range("O508").Value = "Región de Arica y Parinacota"
$ws.Range("P508").Value = 225
$ws.Range("Q508").Value = 80
$ws.Range("R508").Value = "Hortaliza"
